# Update cryptos list prices / 1h volume percentages (GitHub Actions refresh).
# Numeric-looking price strings (e.g. "215.68") are forced to Text format first
# so Excel keeps them as literal strings instead of converting to real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.940.76"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "1.666.20"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("E6").Value = "  +4.66%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0897"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.75%  "
$ws.Range("D12").Value = "1.899.94"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").Value = "1.673.71"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.525"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.65%  "
$ws.Range("D17").Value = "26.927.26"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "236.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.19%  "
$ws.Range("E23").Value = "  -1.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("E32").Value = "  +2.22%  "
$ws.Range("D33").Value = "1.458.68"
$ws.Range("E33").Value = "  -4.63%  "
$ws.Range("E34").Value = "  +2.45%  "
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.905"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0169"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.95%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  +0.56%  "
$ws.Range("E43").Value = "  +6.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "65.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("D45").Value = "1.809.22"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.783"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.87%  "
$ws.Range("E51").Value = "  +0.14%  "
